$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new value, whether to force text (avoid Excel auto-numeric conversion)
$updates = @(
    @{ Cell = "D2"; Value = "67.088.46"; ForceText = $false }
    @{ Cell = "E2"; Value = "  -0.48%  "; ForceText = $false }
    @{ Cell = "D3"; Value = "3.471.95"; ForceText = $false }
    @{ Cell = "E3"; Value = "  -1.59%  "; ForceText = $false }
    @{ Cell = "D5"; Value = "594.28"; ForceText = $true }
    @{ Cell = "E5"; Value = "  -0.42%  "; ForceText = $false }
    @{ Cell = "D6"; Value = "175.65"; ForceText = $true }
    @{ Cell = "E6"; Value = "  +0.94%  "; ForceText = $false }
    @{ Cell = "D7"; Value = "0.999"; ForceText = $true }
    @{ Cell = "E7"; Value = "  -0.02%  "; ForceText = $false }
    @{ Cell = "D8"; Value = "0.586"; ForceText = $true }
    @{ Cell = "E8"; Value = "  -1.22%  "; ForceText = $false }
    @{ Cell = "D9"; Value = "0.129"; ForceText = $true }
    @{ Cell = "E9"; Value = "  -3.99%  "; ForceText = $false }
    @{ Cell = "D10"; Value = "7.07"; ForceText = $true }
    @{ Cell = "E10"; Value = "  -3.21%  "; ForceText = $false }
    @{ Cell = "D11"; Value = "0.423"; ForceText = $true }
    @{ Cell = "E11"; Value = "  -3.17%  "; ForceText = $false }
    @{ Cell = "D12"; Value = "4.064.53"; ForceText = $false }
    @{ Cell = "E12"; Value = "  -1.76%  "; ForceText = $false }
    @{ Cell = "D13"; Value = "31.36"; ForceText = $true }
    @{ Cell = "E13"; Value = "  +8.75%  "; ForceText = $false }
    @{ Cell = "E14"; Value = "  -0.24%  "; ForceText = $false }
    @{ Cell = "D15"; Value = "67.042.14"; ForceText = $false }
    @{ Cell = "E15"; Value = "  -0.39%  "; ForceText = $false }
    @{ Cell = "D16"; Value = "0.0000175"; ForceText = $true }
    @{ Cell = "E16"; Value = "  -4.58%  "; ForceText = $false }
    @{ Cell = "D17"; Value = "3.460.10"; ForceText = $false }
    @{ Cell = "E17"; Value = "  -1.99%  "; ForceText = $false }
    @{ Cell = "D18"; Value = "6.21"; ForceText = $true }
    @{ Cell = "E18"; Value = "  -2.31%  "; ForceText = $false }
    @{ Cell = "D19"; Value = "14.30"; ForceText = $true }
    @{ Cell = "E19"; Value = "  +0.26%  "; ForceText = $false }
    @{ Cell = "D20"; Value = "386.99"; ForceText = $true }
    @{ Cell = "E20"; Value = "  -2.88%  "; ForceText = $false }
    @{ Cell = "D21"; Value = "7.86"; ForceText = $true }
    @{ Cell = "E21"; Value = "  -1.93%  "; ForceText = $false }
    @{ Cell = "E22"; Value = "  +0.20%  "; ForceText = $false }
    @{ Cell = "D23"; Value = "72.64"; ForceText = $true }
    @{ Cell = "E23"; Value = "  -1.25%  "; ForceText = $false }
    @{ Cell = "D24"; Value = "5.71"; ForceText = $true }
    @{ Cell = "E24"; Value = "  -0.08%  "; ForceText = $false }
    @{ Cell = "D25"; Value = "0.532"; ForceText = $true }
    @{ Cell = "E25"; Value = "  -1.64%  "; ForceText = $false }
    @{ Cell = "E26"; Value = "  -1.94%  "; ForceText = $false }
    @{ Cell = "D27"; Value = "10.24"; ForceText = $true }
    @{ Cell = "E27"; Value = "  -0.73%  "; ForceText = $false }
    @{ Cell = "E28"; Value = "  -1.82%  "; ForceText = $false }
    @{ Cell = "E29"; Value = "  -0.27%  "; ForceText = $false }
    @{ Cell = "D30"; Value = "6.08"; ForceText = $true }
    @{ Cell = "E30"; Value = "  -3.60%  "; ForceText = $false }
    @{ Cell = "D31"; Value = "1.41"; ForceText = $true }
    @{ Cell = "E31"; Value = "  -3.87%  "; ForceText = $false }
    @{ Cell = "E32"; Value = "  -2.66%  "; ForceText = $false }
    @{ Cell = "D33"; Value = "23.45"; ForceText = $true }
    @{ Cell = "E33"; Value = "  -2.84%  "; ForceText = $false }
    @{ Cell = "D34"; Value = "7.22"; ForceText = $true }
    @{ Cell = "E34"; Value = "  -2.66%  "; ForceText = $false }
    @{ Cell = "E35"; Value = "  -2.31%  "; ForceText = $false }
    @{ Cell = "D36"; Value = "164.52"; ForceText = $true }
    @{ Cell = "E36"; Value = "  +0.27%  "; ForceText = $false }
    @{ Cell = "D37"; Value = "0.868"; ForceText = $true }
    @{ Cell = "E37"; Value = "  -3.20%  "; ForceText = $false }
    @{ Cell = "E38"; Value = "  -0.54%  "; ForceText = $false }
    @{ Cell = "D39"; Value = "6.89"; ForceText = $true }
    @{ Cell = "E39"; Value = "  -1.14%  "; ForceText = $false }
    @{ Cell = "D40"; Value = "27.06"; ForceText = $true }
    @{ Cell = "E40"; Value = "  -1.53%  "; ForceText = $false }
    @{ Cell = "D41"; Value = "4.59"; ForceText = $true }
    @{ Cell = "E41"; Value = "  -3.20%  "; ForceText = $false }
    @{ Cell = "D42"; Value = "26.06"; ForceText = $true }
    @{ Cell = "E42"; Value = "  -2.02%  "; ForceText = $false }
    @{ Cell = "D43"; Value = "2.788.84"; ForceText = $false }
    @{ Cell = "E43"; Value = "  -0.56%  "; ForceText = $false }
    @{ Cell = "B44"; Value = "dogwifhat"; ForceText = $false }
    @{ Cell = "C44"; Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"; ForceText = $false }
    @{ Cell = "D44"; Value = "2.57"; ForceText = $true }
    @{ Cell = "E44"; Value = "  -2.28%  "; ForceText = $false }
    @{ Cell = "B45"; Value = "Hedera"; ForceText = $false }
    @{ Cell = "C45"; Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"; ForceText = $false }
    @{ Cell = "D45"; Value = "0.0717"; ForceText = $true }
    @{ Cell = "E45"; Value = "  -4.10%  "; ForceText = $false }
    @{ Cell = "D46"; Value = "42.14"; ForceText = $true }
    @{ Cell = "E46"; Value = "  -1.96%  "; ForceText = $false }
    @{ Cell = "E47"; Value = "  -4.70%  "; ForceText = $false }
    @{ Cell = "D48"; Value = "338.81"; ForceText = $true }
    @{ Cell = "E48"; Value = "  -1.06%  "; ForceText = $false }
    @{ Cell = "E49"; Value = "  -3.81%  "; ForceText = $false }
    @{ Cell = "D50"; Value = "33.02"; ForceText = $true }
    @{ Cell = "E50"; Value = "  -2.75%  "; ForceText = $false }
    @{ Cell = "D51"; Value = "6.34"; ForceText = $true }
    @{ Cell = "E51"; Value = "  -3.30%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $u.Value
}
